# Rebuild Sheet2 with the second crawl pass results, matching the canonical diff.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$values = @(
  '{''업체명'': ''한국생명존중희망재단'', ''주소'': ''서울 중구 을지로 6'', ''홈페이지'': ''http://kfsp.org'', ''팩스번호'': ''FAX. (사업총괄본부) 02-3706-0400, (정책지원본부) 02-3706-0401, (자살예방교육부) 02-3706-0402, (자살예방홍보부) '
  '02-3706-0403''}'
  '{''업체명'': ''한국백혈병어린이재단 서울나음소아암센터'', ''주소'': ''서울 성북구 성북로5길 9-14'', ''홈페이지'': ''http://www.kclf.org/'', ''팩스번호'': ''팩스 0507-517-7671''}'
  '{''업체명'': ''한국예술인복지재단'', ''주소'': ''서울 종로구 이화장길 70-15 소호빌딩 1층'', ''홈페이지'': ''http://www.kawf.kr/'', ''팩스번호'': ''''}'
  '{''업체명'': ''롯데재단'', ''주소'': ''서울 중구 남대문로 81 롯데백화점(본점)'', ''홈페이지'': ''http://www.lottefoundation.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''아름다운재단'', ''주소'': ''서울 종로구 자하문로19길 6'', ''홈페이지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''대우재단빌딩'', ''주소'': ''서울 중구 퇴계로 18'', ''홈페이지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''마포복지재단'', ''주소'': ''서울 마포구 대흥로24길 50 1층'', ''홈페이지'': ''http://www.mapowf.or.kr'', ''팩'
  '스번호'': '' ''}'
  '{''업체명'': ''열매나눔재단'', ''주소'': ''서울 중구 퇴계로20길 37'', ''홈페이지'': ''http://www.merryyear.org/'', ''팩스번'
  '호'': ''TEL : 02-310-9508 / FAX : 02-310-9504\nE-mail : csmerry@merryyear.org''}'
  '{''업체명'': ''초록우산 어린이재단'', ''주소'': ''서울 중구 무교로 20 한국복지재단 11층'', ''홈페이지'': ''http://www.childfund.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''종로종합사회복지관'', ''주소'': ''서울 종로구 지봉로13길 82'', ''홈페이지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''계영복지재단'', ''주소'': ''서울 성북구 동소문로 306'', ''홈페이지'': ''http://kehyoung.co.kr/'', ''팩스번호'
  ''': ''팩스: 02-2057-2295''}'
  '{''업체명'': ''대한불교조계종사회복지재단'', ''주소'': ''서울 종로구 우정국로 67'', ''홈페이지'': ''http://jabinanum.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''하나금융나눔재단'', ''주소'': ''서울 중구 다동길 43'', ''홈페이지'': ''http://www.hana-nanum.org/'', ''팩스'
  '번호'': '' ''}'
  '{''업체명'': ''다일복지재단'', ''주소'': ''서울 동대문구 서울시립대로 57 다일복지재단'', ''홈페이지'': '' '', ''팩스번호'': '
  ''' ''}'
  '{''업체명'': ''사회복지법인감리회태화복지재단'', ''주소'': ''서울 종로구 인사동5길 29'', ''홈페이지'': ''http://www.taiwhafound.org/'', ''팩스번호'': ''서울시 종로구 인사동5길 29 태화빌딩 2층 TEL : 02-733-9597 FAX : 02-733-3889''}      '
  '{''업체명'': ''각당복지재단'', ''주소'': ''서울 종로구 경희궁1길 29'', ''홈페이지'': ''http://www.kakdang.or.kr/'', ''팩스'
  '번호'': '' ''}'
  '{''업체명'': ''한국연극인복지재단'', ''주소'': ''서울 종로구 동숭길 29 2층'', ''홈페이지'': ''https://www.facebook.com/playsorkr'', ''팩스번호'': ''전화: 02-741-0332 / 팩스: 02-741-0331''}'
  '{''업체명'': ''송석복지재단'', ''주소'': ''서울 종로구 창경궁로35가길 13'', ''홈페이지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''동아꿈나무재단'', ''주소'': ''서울 서대문구 충정로 29 (주)동아일보사'', ''홈페이지'': ''http://www.dkf.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''재단법인 우체국공익재단'', ''주소'': ''서울 종로구 송월1길 73-10'', ''홈페이지'': ''http://www.kopf.or.kr/'', ''팩스번호'': ''대표전화 : 02-725-3850  팩스 : 02-725-3849 ''}'
  '{''업체명'': ''천태종복지재단'', ''주소'': ''서울 종로구 삼일대로 461 102동, 306호'', ''홈페이지'': ''http://www.with99.org/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''연동복지재단'', ''주소'': ''서울 종로구 대학로3길 9 가나의집 5층 503호 연동복지재단'', ''홈페이지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''사회복지법인성원복지재단'', ''주소'': ''서울 성북구 성북로 87-1 백강빌딩'', ''홈페이지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''KT&G복지재단중부복지센터'', ''주소'': ''서울 중구 서소문로 22'', ''홈페이지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''파라다이스 복지재단'', ''주소'': ''서울 중구 동호로 268 파라다이스빌딩'', ''홈페이지'': ''http://www.paradise.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''IBK행복나눔재단'', ''주소'': ''서울 중구 을지로 79'', ''홈페이지'': ''http://www.ibkfoundation.or.kr/'', '''
  '팩스번호'': ''(재)IBK행복나눔재단 | 서울시 중구 을지로79 ㅣ 대표전화 : 02-3789-3984 ㅣ FAX : 02-3789-3986''}     '
  '{''업체명'': ''한국새생명사랑재단'', ''주소'': ''서울 종로구 율곡로14길 10 오양빌딩4층'', ''홈페이지'': ''http://www.knlf.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''남촌재단'', ''주소'': ''서울 종로구 필운대로 76 GS남촌리더십센터'', ''홈페이지'': ''http://www.jkhuhfoundation.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''보령중보재단'', ''주소'': ''서울 종로구 창경궁로 136 보령빌딩'', ''홈페이지'': ''http://www.boryungjungbo.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''국민사회복지협회'', ''주소'': ''서울 중구 서소문로 89-20 동양빌딩 B동 4층'', ''홈페이지'': ''http://knswa.org/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''일광복지재단'', ''주소'': ''서울 성북구 삼선교로18길 23'', ''홈페이지'': ''http://www.iggywelfare.com/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''평원재단'', ''주소'': ''서울 종로구 혜화로9길 12'', ''홈페이지'': '' '', ''팩스번호'': '' ''}  ''홈페이지'': ''htt'
  '{''업체명'': ''한국장로교복지재단'', ''주소'': ''서울 종로구 대학로3길 29'', ''홈페이지'': ''http://wel.pck.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''진각복지재단'', ''주소'': ''서울 성북구 화랑로13길 17'', ''홈페이지'': ''http://www.jgo.opck.or.kr/'', ''팩r.kr/'', ''팩스번호'': ''서울특별시 성북구 화랑로 13길 17 전화번호 : 02) 942-0144~5 팩스 : 02) 942-0146 E-mail : 029420144@hanmail.net\nCOPYRIGHT (C) 사회복지법인 진각복지재단, ALL RIGHTS REr.kr/'', ''팩스번호SERVED.''}                                                                                    29420144@hanmail.'
  '{''업체명'': ''웅진재단'', ''주소'': ''서울 종로구 새문안로 92 광화문오피시아빌딩 15층'', ''홈페이지'': ''http://www.wjf.kr/'', ''팩스번호'': '' ''}                                                       ''http://www.wjf.'
  '{''업체명'': ''생명문화재단'', ''주소'': ''서울 중구 장충단로6길 15'', ''홈페이지'': '' '', ''팩스번호'': '' ''}                                                                                           ''}'
  '{''업체명'': ''정인욱복지재단'', ''주소'': ''서울 종로구 사직로8길 34'', ''홈페이지'': ''http://www.chungiw-hsf.or.kr/'', giw-hsf.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''골든써클재단'', ''주소'': ''서울 종로구 율곡로 190 여전도회관'', ''홈페이지'': ''http://www.gcfkorea.org/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''국제사랑재단'', ''주소'': ''서울 종로구 대학로 19'', ''홈페이지'': ''http://www.ilovefound.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''양현재단'', ''주소'': ''서울 종로구 북촌로12길 21'', ''홈페이지'': ''http://www.yanghyun.org/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''한국새생명복지재단'', ''주소'': ''서울 성북구 삼선교로24길 55 부영빌딩'', ''홈페이지'': '
  '''http://www.koreassm.org/'', ''팩스번호'': ''팩스 :''}'
  '{''업체명'': ''건양복지'', ''주소'': ''서울 중구 장충단로 226 케이와이 헤리티지 호텔'', ''홈페이지'': ''http://www.keonyang.org/'', ''팩스번호'': ''팩스: 02-902-3333''}'
  '{''업체명'': ''백산복지재단'', ''주소'': ''서울 서대문구 모래내로 285'', ''홈페이지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''사회복지법인닮복지재단'', ''주소'': ''서울 강북구 솔샘로65길 43 그린하우스'', ''홈페이'
  '지'': ''http://www.darmbokji.org/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''어린이재단 서울가정위탁지원센터'', ''주소'': ''서울 중구 무교로 20 어린이재단 빌딩 3'
  '층'', ''홈페이지'': ''http://www.seoul-foster.or.kr/'', ''팩스번호'': ''Fax : 02)325-2664''}'
  '{''업체명'': ''영락사회복지재단'', ''주소'': ''서울 중구 수표로 28'', ''홈페이지'': ''http://www.ynswf.co.kr/'', ''팩스번호'': ''영락사회복지재단      대표자 : 박남진      대표전화 : 02-2265-7023      '
  '팩스 : 02-2277-5855\n주소 : 서울시 중구 수표로 28 보아스 501호      사업자등록번호 : 202-82-30130      이메일 : ynswf@naver.com''}'
  '{''업체명'': ''세종마을푸르메센터'', ''주소'': ''서울 종로구 자하문로 89'', ''홈페이지'': ''http://www.purme.org/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''서울시50플러스 서부캠퍼스'', ''주소'': ''서울 은평구 통일로62길 7'', ''홈페이지'': '' '', '
  '''팩스번호'': '' ''}'
  '{''업체명'': ''다일공동체 법인사무국'', ''주소'': ''서울 동대문구 서울시립대로 57'', ''홈페이지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''KT&G복지재단 북부복지센터'', ''주소'': ''서울 강북구 수유로 50 젬스톤타워'', ''홈페이지'
  ''': ''http://www.ktngwelfare.org/welfare/welfare_02.asp'', ''팩스번호'': '' ''}'
  '{''업체명'': ''여전도회작은자복지재단'', ''주소'': ''서울 종로구 율곡로 190'', ''홈페이지'': '' '', ''팩스'
  '번호'': '' ''}'
  '{''업체명'': ''서울복지재단 장애인전환서비스지원센터'', ''주소'': ''서울 서대문구 통일로 135 충정빌'
  '딩'', ''홈페이지'': ''http://happy.welfare.seoul.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''오병이어복지재단 럽앤홉'', ''주소'': ''서울 은평구 역말로10길 30-1'', ''홈페이지'': ''http://www.xn--dm2b48vz8m.org/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''김남호복지재단'', ''주소'': ''서울 중구 명동길 74'', ''홈페이지'': '' '', ''팩스번호'': '' ''}{''업체명'': ''콘라드 아데나워 재단'', ''주소'': ''서울 종로구 창덕궁1길 13 원서빌딩 2층'', ''홈페이지'
  ''': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''우리모두복지재단 선하어린이집'', ''주소'': ''서울 서대문구 북아현로4자길 31'', ''홈페이'
  '지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''구세군복지재단홍은종합사회복지관홍제동별관'', ''주소'': ''서울 서대문구 세검정로 71'', ''홈페이지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''다정한사람들복지재단'', ''주소'': ''서울 종로구 종로9길 8'', ''홈페이지'': '' '', ''팩스번'
  '호'': '' ''}'
  '{''업체명'': ''선암사회복지재단'', ''주소'': ''서울 종로구 성균관로 65'', ''홈페이지'': ''http://www.sunam.net/'', ''팩스번호'': ''주소 (03070) 서울특별시 종로구 성균관로 65 (구.서울특별시 종로구 명륜 '
  '동1가 10-1) 고객센터:02-745-1932 팩스:02-744-0496\n사업자 번호 : 208-82-03795 상호 : 선암사회'
  '복지재단 대표 : 황갑성''}'
  '{''업체명'': ''사회복지법인애원복지재단'', ''주소'': ''서울 중구 서애로 17'', ''홈페이지'': '' '', ''팩스'
  '번호'': '' ''}'
  '{''업체명'': ''사단법인 장아람재단'', ''주소'': ''서울 서대문구 신촌로 129 아트레온'', ''홈페이지'': ''https://blog.naver.com/jangaram1995'', ''팩스번호'': '' ''}'
  '{''업체명'': ''유엔한국평화복지재단'', ''주소'': ''서울 성북구 동소문로 69 동문빌딩'', ''홈페이지'': ''http://www.unipwa.net/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''우체국한사랑의집'', ''주소'': ''서울 종로구 대학로5길 27 우체국한사랑의집'', ''홈페이지'
  ''': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''환주복지재단'', ''주소'': ''서울 중구 명동길 74'', ''홈페이지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''기쁨나눔재단'', ''주소'': ''서울 마포구 백범로5길 13 창선빌딩 B동 1층'', ''홈페이지'': ''https://www.instagram.com/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''한국소아암재단'', ''주소'': ''서울 중구 중림로 47-1'', ''홈페이지'': ''http://www.angelc.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''참좋은나눔재단'', ''주소'': ''서울 서대문구 통일로 371 서강빌딩 4층 참좋은나눔재단'', '
  '''홈페이지'': ''http://charmnanum.org'', ''팩스번호'': '' ''}'
  '{''업체명'': ''노인의료나눔재단'', ''주소'': ''서울 중구 통일로 86 바비엥3차 7F 701호'', ''홈페이지'': '
  '''http://www.ok6595.or.kr/'', ''팩스번호'': ''FAX''}'
  '{''업체명'': ''사회복지법인어린이재단'', ''주소'': ''서울 중구 무교로 20'', ''홈페이지'': '' '', ''팩스번'
  '호'': '' ''}'
  '{''업체명'': ''무지개청소년센터'', ''주소'': ''서울 종로구 자하문로24길 20'', ''홈페이지'': '' '', ''팩스'
  '번호'': '' ''}'
  '{''업체명'': ''재단법인 동천'', ''주소'': ''서울 종로구 우정국로 26 센트로폴리스빌딩 b동 10층'', ''홈'
  '페이지'': ''http://www.bkl.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''한마음한몸운동본부'', ''주소'': ''서울 중구 명동길 66'', ''홈페이지'': ''http://www.obos.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''대산문화재단'', ''주소'': ''서울 종로구 종로 1 교보생명빌딩'', ''홈페이지'': ''http://www.daesan.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''아름다운동행'', ''주소'': ''서울 종로구 인사동11길 22 다연 2층'', ''홈페이지'': ''http://www.dreaminus.org/'', ''팩스번호'': ''이용약관\n개인정보처리방침\n이메일주소무단수집거부\n문의하 '
  '기\n저작권 정책\n국세청\n행정안전부\n재단법인 아름다운동행 03145 서울시 종로구 우정국로 55 ( '
  '견지동)\n대표자 이규정   사업자등록번호 101-82-17154\n전화 02) 737-9595   팩스 02) 737-9195  '
  ' 이메일 thenanum@hanmail.net''}'
  '{''업체명'': ''한국의료지원재단'', ''주소'': ''서울 중구 남대문로10길 9 경기빌딩'', ''홈페이지'': ''http://komaf12.org/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''사회복지법인 한국노인생활지원재단'', ''주소'': ''서울 동대문구 천호대로 17'', ''홈페이'
  '지'': ''http://www.hid081.co.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''한국재활재단'', ''주소'': ''서울 서대문구 연희로 248'', ''홈페이지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''어린이재단 서울지역본부'', ''주소'': ''서울 중구 무교로 20 어린이재단빌딩 3층'', ''홈페'
  '이지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''함께나누는세상'', ''주소'': ''서울 종로구 통일로 162 덕산빌딩 201호'', ''홈페이지'': ''http://www.sharingtogether.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''재단법인 신용카드사회공헌재단'', ''주소'': ''서울 중구 다동길 43 한외빌딩 11층'', ''홈'
  '페이지'': ''http://www.ccfd.or.kr/main/main.html'', ''팩스번호'': ''팩스 : 02-2135-1837''}'
  '{''업체명'': ''동숭학술재단'', ''주소'': ''서울 종로구 동숭3길 27'', ''홈페이지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''내고장사랑재단'', ''주소'': ''서울 종로구 율곡로2길 7 로얄팰리스'', ''홈페이지'': ''http://www.myhomelove.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''베스티안재단 사회복지사업본부'', ''주소'': ''서울 성동구 왕십리로 382 1층'', ''홈페이지'
  ''': ''http://ibestian.org/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''에이드투더처치인니드코리아'', ''주소'': ''서울 중구 명동길 66'', ''홈페이지'': ''http://www.churchinneed.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''기독교청년회유지재단'', ''주소'': ''서울 종로구 종로 69'', ''홈페이지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''세계평화나눔재단'', ''주소'': ''서울 종로구 새문안로 92 광화문오피시아빌딩'', ''홈페이'
  '지'': ''http://unpsf.co.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''국제밝은사회재단'', ''주소'': ''서울 종로구 성균관로15가길 38'', ''홈페이지'': '' '', ''팩'
  '스번호'': '' ''}'
  '{''업체명'': ''삼선배움과나눔재단'', ''주소'': ''서울 중구 퇴계로 36 삼선빌딩'', ''홈페이지'': ''http://www.sscare.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''사단법인 웨슬리나눔재단'', ''주소'': ''서울 종로구 세종대로 149 16층'', ''홈페이지'': '' '
  ''', ''팩스번호'': '' ''}'
  '{''업체명'': ''국제나눔재단'', ''주소'': ''서울 종로구 삼일대로 469 서원빌딩 15층'', ''홈페이지'': ''https://innf.co.kr/'', ''팩스번호'': ''Tel. 02-723-5400 | Fax. 02-735-4500 | innf21@daum.net\n주소. '
  '서울시 종로구 삼일대로 469, 서원빌딩 15층 ㅣ 고유번호. 306-82-69048''}'
  '{''업체명'': ''어린이재단 서울남부지역본부'', ''주소'': ''서울 중구 무교로 20 어린이재단'', ''홈페이지'
  ''': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''국제나눔재단'', ''주소'': ''서울 성북구 창경궁로 325'', ''홈페이지'': ''http://internationalnanum.com/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''재단법인 통일아동키크기재단'', ''주소'': ''서울 마포구 와우산로 121'', ''홈페이지'': ''http://www.unigrow.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''재단법인 동국장학회'', ''주소'': ''서울 종로구 인사동7길 12 백상빌딩 302호'', ''홈페이'
  '지'': ''http://dongguk.or.kr'', ''팩스번호'': '' ''}'
  '{''업체명'': ''진실의 힘'', ''주소'': ''서울 중구 세종대로19길 16'', ''홈페이지'': ''http://www.truthfoundation.or.kr/'', ''팩스번호'': ''Fax''}'
  '{''업체명'': ''서비스포피스'', ''주소'': ''서울 서대문구 성산로 527 하늬솔빌딩'', ''홈페이지'': ''http://www.sfp.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''사회복지법인인덕원'', ''주소'': ''서울 은평구 서오릉로 30-17'', ''홈페이지'': '' '', ''팩스'
  '번호'': '' ''}'
  '{''업체명'': ''한음문화재단'', ''주소'': ''서울 서대문구 충정로 76'', ''홈페이지'': ''http://www.hecf.or.kr/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''평화의 길 국제재단'', ''주소'': ''서울 종로구 새문안로 42'', ''홈페이지'': ''http://road1004.onmam.com/'', ''팩스번호'': '' ''}'
  '{''업체명'': ''식량나눔재단'', ''주소'': ''서울 서대문구 이화여대5길 15'', ''홈페이지'': ''http://foodnanum.org/'', ''팩스번호'': ''서울특별시 중구 을지로 218, 421-2호 / FAX 050-2393-3333''}'
  '{''업체명'': ''동승문화재단'', ''주소'': ''서울 종로구 청계천로 279'', ''홈페이지'': '' '', ''팩스번호'': '' ''}'
  '{''업체명'': ''국내장학문화재단'', ''주소'': ''서울 종로구 종로1길 36'', ''홈페이지'': '' '', ''팩스번호'': '' ''}'
)

# A scratch cell well outside the used range, for the apostrophe-preserving paste trick.
$scratch = $ws2.Cells.Item(2000, 200)

function Set-CellText($cell, [string]$text) {
  if ($text.Length -gt 0 -and $text.Substring(0,1) -eq "'") {
    # A literal leading apostrophe would otherwise be swallowed as a
    # quote-prefix marker by the normal Value setter, so build it via a
    # formula and paste the computed value back in as a literal.
    $script:scratch.Formula = '="' + $text.Substring(0,1) + '"' + '&"' + $text.Substring(1) + '"'
    $script:scratch.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
  } else {
    $cell.Value = $text
  }
}

# First pass: populate every row except 4 and 8 in order.
for ($i = 0; $i -lt $values.Count; $i++) {
  $r = $i + 1
  if ($r -eq 4 -or $r -eq 8) { continue }
  Set-CellText $ws2.Cells.Item($r, 1) $values[$i]
}

# Error-check fixup pass: rows 4 then 8 get (re)written last.
Set-CellText $ws2.Cells.Item(4, 1) $values[3]
Set-CellText $ws2.Cells.Item(8, 1) $values[7]

$scratch.Clear() | Out-Null
$excel.CutCopyMode = $false

# Activate Sheet2 (making it the selected/visible tab) and set its selection.
$ws2.Activate()
$ws2.Range("K6").Select() | Out-Null
